$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Spouse" column (E) with pseudonymized spouse names
$ws.Range("E1").Value = "Spouse"
$ws.Range("E2").Value = "Peggy Hill"
$ws.Range("E5").Value = "Marge Simpson"

# Set column E width to match bestFit sizing used for the other columns
$ws.Columns.Item(5).ColumnWidth = 14.7109375

# Select the new column, matching the author's selection after adding it
$ws.Columns.Item(5).Select()
